$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.260887446824029
$ws.Cells.Item(2, 3).Value = 0.199249073379633
$ws.Cells.Item(2, 4).Value = 0.07823181986877614
$ws.Cells.Item(2, 5).Value = 0.09758177810482982
$ws.Cells.Item(2, 7).Value = 0.002463313228403287
$ws.Cells.Item(2, 12).Value = 0.1899887431010754
$ws.Cells.Item(2, 13).Value = 0.2663511255113491
$ws.Cells.Item(2, 14).Value = 1.55301995052443
$ws.Cells.Item(2, 15).Value = 4.054209819659405
$ws.Cells.Item(3, 2).Value = 1.171044405421299
$ws.Cells.Item(3, 3).Value = 0.1854774821674425
$ws.Cells.Item(3, 4).Value = 0.07099981640112674
$ws.Cells.Item(3, 5).Value = 0.09825067771945184
$ws.Cells.Item(3, 7).Value = 0.00246683245833057
$ws.Cells.Item(3, 12).Value = 0.1874107202957163
$ws.Cells.Item(3, 13).Value = 0.2525465544461696
$ws.Cells.Item(3, 14).Value = 1.571932313769199
$ws.Cells.Item(3, 15).Value = 4.03968725914379
$ws.Cells.Item(4, 2).Value = 1.116318715562272
$ws.Cells.Item(4, 3).Value = 0.1769469899430902
$ws.Cells.Item(4, 4).Value = 0.06659608434817699
$ws.Cells.Item(4, 5).Value = 0.09868527736662824
$ws.Cells.Item(4, 7).Value = 0.002469109539887057
$ws.Cells.Item(4, 12).Value = 0.1859212253667764
$ws.Cells.Item(4, 13).Value = 0.2441849015072748
$ws.Cells.Item(4, 14).Value = 1.584142274297843
$ws.Cells.Item(4, 15).Value = 4.033347244086116
$ws.Cells.Item(5, 2).Value = 1.09412854058948
$ws.Cells.Item(5, 3).Value = 0.1734520390257757
$ws.Cells.Item(5, 4).Value = 0.06481071855714049
$ws.Cells.Item(5, 5).Value = 0.09886839737826003
$ws.Cells.Item(5, 7).Value = 0.002470066794798043
$ws.Cells.Item(5, 12).Value = 0.1853377690616043
$ws.Cells.Item(5, 13).Value = 0.2408063511415222
$ws.Cells.Item(5, 14).Value = 1.589268168265596
$ws.Cells.Item(5, 15).Value = 4.031410647027315
$ws.Cells.Item(6, 2).Value = 1.090450605305278
$ws.Cells.Item(6, 3).Value = 0.1728705774511639
$ws.Cells.Item(6, 4).Value = 0.06451481358588751
$ws.Cells.Item(6, 5).Value = 0.09889916800542209
$ws.Cells.Item(6, 7).Value = 0.002470227520011855
$ws.Cells.Item(6, 12).Value = 0.1852423084303254
$ws.Cells.Item(6, 13).Value = 0.2402470939986472
$ws.Cells.Item(6, 14).Value = 1.590128391285493
$ws.Cells.Item(6, 15).Value = 4.031128126172547
$ws.Cells.Item(7, 2).Value = 1.116019000688254
$ws.Cells.Item(7, 3).Value = 0.176899931455651
$ws.Cells.Item(7, 4).Value = 0.06657196913280927
$ws.Cells.Item(7, 5).Value = 0.09868772261168868
$ws.Cells.Item(7, 7).Value = 0.002469122331048571
$ws.Cells.Item(7, 12).Value = 0.1859132613743668
$ws.Cells.Item(7, 13).Value = 0.2441392200543504
$ws.Cells.Item(7, 14).Value = 1.584210795677532
$ws.Cells.Item(7, 15).Value = 4.033318508038235
$ws.Cells.Item(8, 2).Value = 1.229819029629994
$ws.Cells.Item(8, 3).Value = 0.1945161681263698
$ws.Cells.Item(8, 4).Value = 0.07573056603993678
$ws.Cells.Item(8, 5).Value = 0.09780746133108598
$ws.Cells.Item(8, 7).Value = 0.002464502583775546
$ws.Cells.Item(8, 12).Value = 0.1890804698584816
$ws.Cells.Item(8, 13).Value = 0.2615676317610536
$ws.Cells.Item(8, 14).Value = 1.559416848011971
$ws.Cells.Item(8, 15).Value = 4.048667003981052
$ws.Cells.Item(9, 2).Value = 1.456436650936894
$ws.Cells.Item(9, 3).Value = 0.2284673569338054
$ws.Cells.Item(9, 4).Value = 0.09398531071138905
$ws.Cells.Item(9, 5).Value = 0.09627042222974502
$ws.Cells.Item(9, 7).Value = 0.002456361598726715
$ws.Cells.Item(9, 12).Value = 0.196031707576438
$ws.Cells.Item(9, 13).Value = 0.2966491485530725
$ws.Cells.Item(9, 14).Value = 1.51554039127795
$ws.Cells.Item(9, 15).Value = 4.099263395914733
$ws.Cells.Item(10, 2).Value = 1.625028815287237
$ws.Cells.Item(10, 3).Value = 0.2530491241708717
$ws.Cells.Item(10, 4).Value = 0.1075825549869762
$ws.Cells.Item(10, 5).Value = 0.09525586136941655
$ws.Cells.Item(10, 7).Value = 0.00245093444353396
$ws.Cells.Item(10, 12).Value = 0.2015897316295252
$ws.Cells.Item(10, 13).Value = 0.3229734278406227
$ws.Cells.Item(10, 14).Value = 1.486198424283351
$ws.Cells.Item(10, 15).Value = 4.149015209203526
$ws.Cells.Item(11, 2).Value = 1.702179869630299
$ws.Cells.Item(11, 3).Value = 0.2641535249395588
$ws.Cells.Item(11, 4).Value = 0.1138098922097441
$ws.Cells.Item(11, 5).Value = 0.09481908779182358
$ws.Cells.Item(11, 7).Value = 0.00244858455547759
$ws.Cells.Item(11, 12).Value = 0.2042161171225132
$ws.Cells.Item(11, 13).Value = 0.3350683127249425
$ws.Cells.Item(11, 14).Value = 1.473478846064667
$ws.Cells.Item(11, 15).Value = 4.174398501286362
$ws.Cells.Item(12, 2).Value = 1.73146027963287
$ws.Cells.Item(12, 3).Value = 0.2683472207487512
$ws.Cells.Item(12, 4).Value = 0.1161741208323122
$ws.Cells.Item(12, 5).Value = 0.09465724342046755
$ws.Cells.Item(12, 7).Value = 0.002447711725161265
$ws.Cells.Item(12, 12).Value = 0.2052247362764774
$ws.Cells.Item(12, 13).Value = 0.3396654884962942
$ws.Cells.Item(12, 14).Value = 1.468752694138894
$ws.Cells.Item(12, 15).Value = 4.184407380262144
$ws.Cells.Item(13, 2).Value = 1.725151342544564
$ws.Cells.Item(13, 3).Value = 0.2674445372769583
$ws.Cells.Item(13, 4).Value = 0.115664671109343
$ws.Cells.Item(13, 5).Value = 0.0946919416303218
$ws.Cells.Item(13, 7).Value = 0.002447898948999769
$ws.Cells.Item(13, 12).Value = 0.2050068870722157
$ws.Cells.Item(13, 13).Value = 0.3386746456904675
$ws.Cells.Item(13, 14).Value = 1.469766529315587
$ws.Cells.Item(13, 15).Value = 4.182234123224362
$ws.Cells.Item(14, 2).Value = 1.704587489002733
$ws.Cells.Item(14, 3).Value = 0.2644987697826764
$ws.Cells.Item(14, 4).Value = 0.1140042767693359
$ws.Cells.Item(14, 5).Value = 0.094805701607098
$ws.Cells.Item(14, 7).Value = 0.002448512406259177
$ws.Cells.Item(14, 12).Value = 0.204298815135914
$ws.Cells.Item(14, 13).Value = 0.3354461830510616
$ws.Cells.Item(14, 14).Value = 1.473088208616181
$ws.Cells.Item(14, 15).Value = 4.175213978912495
$ws.Cells.Item(15, 2).Value = 1.691999972819133
$ws.Cells.Item(15, 3).Value = 0.2626929279707326
$ws.Cells.Item(15, 4).Value = 0.1129880293750603
$ws.Cells.Item(15, 5).Value = 0.09487584527189563
$ws.Cells.Item(15, 7).Value = 0.002448890381107657
$ws.Cells.Item(15, 12).Value = 0.2038669312850914
$ws.Cells.Item(15, 13).Value = 0.3334708803267503
$ws.Cells.Item(15, 14).Value = 1.475134619950833
$ws.Cells.Item(15, 15).Value = 4.170965645484443
$ws.Cells.Item(16, 2).Value = 1.619995950375085
$ws.Cells.Item(16, 3).Value = 0.2523218490057957
$ws.Cells.Item(16, 4).Value = 0.1071764307887264
$ws.Cells.Item(16, 5).Value = 0.09528490304871107
$ws.Cells.Item(16, 7).Value = 0.002451090400619111
$ws.Cells.Item(16, 12).Value = 0.2014200608413574
$ws.Cells.Item(16, 13).Value = 0.3221853992232369
$ws.Cells.Item(16, 14).Value = 1.487042331874466
$ws.Cells.Item(16, 15).Value = 4.147411809376251
$ws.Cells.Item(17, 2).Value = 1.575940435508528
$ws.Cells.Item(17, 3).Value = 0.2459394956053416
$ws.Cells.Item(17, 4).Value = 0.1036219649874681
$ws.Cells.Item(17, 5).Value = 0.0955421818205302
$ws.Cells.Item(17, 7).Value = 0.002452470448357173
$ws.Cells.Item(17, 12).Value = 0.199944066902205
$ws.Cells.Item(17, 13).Value = 0.3152927247277546
$ws.Cells.Item(17, 14).Value = 1.494508328531312
$ws.Cells.Item(17, 15).Value = 4.13366772457249
$ws.Cells.Item(18, 2).Value = 1.550643992265918
$ws.Cells.Item(18, 3).Value = 0.2422612012081231
$ws.Cells.Item(18, 4).Value = 0.1015814695743558
$ws.Cells.Item(18, 5).Value = 0.09569249228930587
$ws.Cells.Item(18, 7).Value = 0.002453275416381682
$ws.Cells.Item(18, 12).Value = 0.1991043414068514
$ws.Cells.Item(18, 13).Value = 0.3113395279494782
$ws.Cells.Item(18, 14).Value = 1.498861695737983
$ws.Cells.Item(18, 15).Value = 4.126021350388186
$ws.Cells.Item(19, 2).Value = 1.542086483062178
$ws.Cells.Item(19, 3).Value = 0.2410145370745056
$ws.Cells.Item(19, 4).Value = 0.100891267986114
$ws.Cells.Item(19, 5).Value = 0.0957437853296117
$ws.Cells.Item(19, 7).Value = 0.002453549891099338
$ws.Cells.Item(19, 12).Value = 0.198821610481275
$ws.Cells.Item(19, 13).Value = 0.3100029859398248
$ws.Cells.Item(19, 14).Value = 1.500345821467249
$ws.Cells.Item(19, 15).Value = 4.123476842616299
$ws.Cells.Item(20, 2).Value = 1.580625763142677
$ws.Cells.Item(20, 3).Value = 0.2466196676535901
$ws.Cells.Item(20, 4).Value = 0.1039999358866197
$ws.Cells.Item(20, 5).Value = 0.09551455291013311
$ws.Cells.Item(20, 7).Value = 0.002452322381225757
$ws.Cells.Item(20, 12).Value = 0.200100234217814
$ws.Cells.Item(20, 13).Value = 0.3160252944273765
$ws.Cells.Item(20, 14).Value = 1.493707440404732
$ws.Cells.Item(20, 15).Value = 4.135104006185315
$ws.Cells.Item(21, 2).Value = 1.710625834284542
$ws.Cells.Item(21, 3).Value = 0.2653643204725711
$ws.Cells.Item(21, 4).Value = 0.114491809680132
$ws.Cells.Item(21, 5).Value = 0.09477219119510227
$ws.Cells.Item(21, 7).Value = 0.002448331757891175
$ws.Cells.Item(21, 12).Value = 0.2045064115737318
$ws.Cells.Item(21, 13).Value = 0.3363939975622827
$ws.Cells.Item(21, 14).Value = 1.472110094264982
$ws.Cells.Item(21, 15).Value = 4.177265187495038
$ws.Cells.Item(22, 2).Value = 1.79596676983715
$ws.Cells.Item(22, 3).Value = 0.2775491305636137
$ws.Cells.Item(22, 4).Value = 0.1213842655303381
$ws.Cells.Item(22, 5).Value = 0.09430771702186114
$ws.Cells.Item(22, 7).Value = 0.002445822830617291
$ws.Cells.Item(22, 12).Value = 0.2074680650088254
$ws.Cells.Item(22, 13).Value = 0.3498057672668935
$ws.Cells.Item(22, 14).Value = 1.458522412477626
$ws.Cells.Item(22, 15).Value = 4.207133236628238
$ws.Cells.Item(23, 2).Value = 1.750384356980589
$ws.Cells.Item(23, 3).Value = 0.2710519272619365
$ws.Cells.Item(23, 4).Value = 0.1177023761372595
$ws.Cells.Item(23, 5).Value = 0.09455372383504068
$ws.Cells.Item(23, 7).Value = 0.002447152845157596
$ws.Cells.Item(23, 12).Value = 0.2058798851844728
$ws.Cells.Item(23, 13).Value = 0.3426385776663352
$ws.Cells.Item(23, 14).Value = 1.465726107136922
$ws.Cells.Item(23, 15).Value = 4.190980035009318
$ws.Cells.Item(24, 2).Value = 1.578507427840805
$ws.Cells.Item(24, 3).Value = 0.2463121898841791
$ws.Cells.Item(24, 4).Value = 0.1038290458600954
$ws.Cells.Item(24, 5).Value = 0.09552703647222449
$ws.Cells.Item(24, 7).Value = 0.002452389286384384
$ws.Cells.Item(24, 12).Value = 0.2000296034277795
$ws.Cells.Item(24, 13).Value = 0.3156940700852147
$ws.Cells.Item(24, 14).Value = 1.49406933170947
$ws.Cells.Item(24, 15).Value = 4.134453868144021
$ws.Cells.Item(25, 2).Value = 1.394761711380909
$ws.Cells.Item(25, 3).Value = 0.2193462688196917
$ws.Cells.Item(25, 4).Value = 0.08901476338748182
$ws.Cells.Item(25, 5).Value = 0.09666604541710955
$ws.Cells.Item(25, 7).Value = 0.002463313228403287
$ws.Cells.Item(25, 12).Value = 0.1940719952863859
$ws.Cells.Item(25, 13).Value = 0.2870619983918061
$ws.Cells.Item(25, 14).Value = 1.526902052738278
$ws.Cells.Item(25, 15).Value = 4.083372772845735
